$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362, shifting existing rows 362-453 down to 363-454
$ws.Range("A362").EntireRow.Insert()

# Populate the newly inserted row 362 with its data.
# Columns A, B, C, E, F, G, I, R carry over the same values as the row that used
# to occupy position 362 (now shifted to 363), columns D, H, J, K, L, M, N, O, P, Q
# hold new values.
$ws.Range("A362").Value2 = 9
$ws.Range("B362").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C362").Value2 = "Metropolitana"
$ws.Range("D362").Value2 = 45209
$ws.Range("E362").Value2 = 13
$ws.Range("F362").Value2 = 100112021
$ws.Range("G362").Value2 = "Ají"
$ws.Range("H362").Value2 = "Inferno"
$ws.Range("I362").Value2 = "Primera"
$ws.Range("J362").Value2 = 34
$ws.Range("K362").Value2 = 29000
$ws.Range("L362").Value2 = 30000
$ws.Range("M362").Value2 = 29500
$ws.Range("N362").Value2 = "$/caja 10 kilos"
$ws.Range("O362").Value2 = "Región de Arica y Parinacota"
$ws.Range("P362").Value2 = 2950
$ws.Range("Q362").Value2 = 10
$ws.Range("R362").Value2 = "Hortaliza"
